$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the pT bin index in column B (rows 2-17) by 1, shifting the
# "pT < 1 GeV" bin (formerly row 18, index 16) out of the table.
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    $cell.Value = $current + 1
}

# Clear out the last row (row 18), which held the removed pT < 1 GeV bin.
$ws.Range("A18:J18").ClearContents()

# Move the active selection to B19, matching the author's post-edit cursor.
$ws.Range("B19").Select()
